$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D and E columns are treated as text so numeric-looking / percent strings are preserved exactly
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.263.95"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.16%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.853.36"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.21%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.66"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.53%  "

$ws.Range("B6").Value = "XRP"
$ws.Range("C6").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6917"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -6.75%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.11%  "

$ws.Range("B8").Value = "Dogecoin"
$ws.Range("C8").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07689"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +6.65%  "

$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3060"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.83%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.59"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -4.37%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08092"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.34%  "

$ws.Range("B12").Value = "Polygon"
$ws.Range("C12").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7231"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.88%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.830.33"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -5.85%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.207"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.84%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.09"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.76%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.299.25"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.21%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.881"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.15%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "241.58"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -4.09%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.75%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.11"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.46%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.0000"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.11%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.101.77"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.55%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.02%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.607"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -5.57%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.051"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.20%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.11"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.75%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1455"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -6.14%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.08"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.47%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -5.06%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.398"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -7.46%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.509"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.85%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.422"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.99%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.042"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -5.82%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05223"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.89%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.195"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.39%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7140"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -4.60%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9963"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.81%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.671"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.78%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01859"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -5.27%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.698"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.18%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8936"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +4.52%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4295"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -5.45%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.860"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.90%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "69.65"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.94%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.043.64"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -6.17%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.02%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.40"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.41%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.260"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -4.56%  "

$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.736"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -6.38%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.998.24"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.73%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.246"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.39%  "
